# Complete test for Specific_Risk macro
#
# The test workbook lists, one row per macro test, the macro name (A),
# its description (B) and the test identifier (C). The "Specific Risk"
# test row (row 13) gets a more descriptive label, and a brand-new test
# row is appended at the bottom of the table for the scale=252 variant
# of the Specific_Risk macro.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 / column B: flesh out the description of the existing
# Table_SpecificRisk test to mention the scale used.
$ws.Range("A64").Value = "Specific_Risk"
$ws.Range("B13").Value = "Test Specific Risk table with scale=252"
$ws.Range("B64").Value = "Test Specific Risk with scale=252"
$ws.Range("C64").Value = "Specific_Risk_test"

# Leave the selection on the newly-added last cell, matching how the
# workbook was left after the edit.
$ws.Range("C64").Select()
